# articles_Piper_Jan_1.xlsx manual data update
# - removes the now-unused "datetime" and "img" columns
# - shifts "desc", "link" (hyperlink) and "summary" columns left accordingly
# - re-applies the Hyperlink style/relationship to the (now relocated) link column
# - lightly reformats row 3 (mirrors a "reselect & restyle" edit made by the author)
# - tidies up the sheet view (zoom, selection) and column widths

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- known hyperlink targets (read before we touch anything) ---
$link1 = "https://www.nzz.ch/feuilleton/annemarie-pieper-die-philosophin-gibt-mit-80-neue-denkanstoesse-ld.1594828"
$link2 = "https://www.ndr.de/kultur/buch/Lionel-Shriver-Die-perfekte-Freundin,perfektefreundin102.html"

# --- remove the "datetime" column (was column E, always empty) ---
$ws.Columns("E").Delete()

# --- remove the "img" column (was column H, now column G after the delete above) ---
$ws.Columns("G").Delete()

# --- the hyperlinks collection does not automatically follow the column shift,
#     so rebuild it pointing at the new location of the "link" column (F) ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), $link1)
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), $link2)
$ws.Range("F3").Style = "Hyperlink"

# --- minor reformat touch on row 3 (no visible color change, explicit reapply) ---
$ws.Range("A3:G3").Interior.ColorIndex = -4142

# --- sheet view: zoom + select row 3 like in the saved workbook ---
$ws.Application.ActiveWindow.Zoom = 71
$ws.Rows("3").Select()

# --- column widths ---
$ws.Columns("D").ColumnWidth = 9.14
$ws.Columns("E").ColumnWidth = 8.43
$ws.Columns("F").ColumnWidth = 8.43

# --- page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "applied edits"
